$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column B (shifts existing B:M to D:O)
$ws.Range("B:C").Insert()

# New header cells (row 1), matching the style of the existing bold/centered header row
$ws.Range("B1").Value = "Count"
$ws.Range("C1").Value = "Patients"
$ws.Range("B1:C1").Font.Bold = $true
$ws.Range("B1:C1").HorizontalAlignment = -4108

# New data values: Count column (B) and Patients column (C)
$counts = @(2999, 2999, 2989, 2999, 2989, 2905, 2905, 1689, 2908)
for ($i = 0; $i -lt $counts.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $counts[$i]
    $ws.Cells.Item($row, 3).Value = 380
}
